# "Replace old MQTT timeout timer code to prevent stall condition.
#  Introduced mqttRuntime function as solution."
#
# The symbol table on Tabelle1 (sheet1) gets a new row for the freshly
# introduced mqttRuntime FB (FB 102). It is inserted right after the
# existing "MQTT" / "FB 100" / mqttPacketReader / "FB 101" block (row 7),
# pushing everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The sheet is protected (format-only protection, cells unlocked) - lift it
# so the rows can be edited; the target workbook no longer carries
# sheetProtection at all, so we leave it unprotected.
$ws.Unprotect("840F")

# Insert a new row at position 7; existing rows 7-21 shift down to 8-22.
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "mqttRuntime"
$ws.Range("B7").Value = "FB      102"
$ws.Range("C7").Value = "FB      102"

# The row-insert carries the (empty, but styled) column D cell all the way
# down to the new last row (D22), but the authored workbook's final row
# has no column D cell at all - match that by resetting D22 back to the
# plain default style and clearing it, which drops the cell entirely.
$ws.Range("D22").Style = $ws.Range("A22").Style
$ws.Range("D22").ClearContents()

$ws.Range("C8").Select()
